$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5047.975
$ws.Range("J43").Value = 4612.4
$ws.Range("L43").Value = 4612.4
$ws.Range("N43").Value = -4750.4

$ws.Range("H53").Value = 1692.875
$ws.Range("I53").Value = 2170.75
$ws.Range("J53").Value = 1215
$ws.Range("K53").Value = 2170.75
$ws.Range("L53").Value = 1215
$ws.Range("M53").Value = -1533.75
$ws.Range("N53").Value = -2489

$ws.Range("H64").Value = 6259.8
$ws.Range("J64").Value = 5674.625
$ws.Range("L64").Value = 5674.625
$ws.Range("N64").Value = -6170.625

$ws.Range("H67").Value = 6259.8
$ws.Range("J67").Value = 5674.625
$ws.Range("L67").Value = 5674.625
$ws.Range("N67").Value = -7390.625

$ws.Range("H99").Value = 3141.0715
$ws.Range("I99").Value = 2490.75
$ws.Range("J99").Value = 4008.1667
$ws.Range("K99").Value = 7472.25
$ws.Range("L99").Value = 12024.5001
$ws.Range("M99").Value = -5974.25
$ws.Range("N99").Value = -15020.5001

$ws.Range("H116").Value = 47729.043
$ws.Range("I116").Value = 75980.07000000001
$ws.Range("K116").Value = 75980.07000000001
$ws.Range("M116").Value = -72538.07000000001

$ws.Range("H132").Value = 103386.6
$ws.Range("I132").Value = 128881.875
$ws.Range("K132").Value = 386645.625
$ws.Range("M132").Value = -384115.625

$ws.Range("H138").Value = 6696.912
$ws.Range("I138").Value = 3494.4375
$ws.Range("J138").Value = 7946.6587
$ws.Range("K138").Value = 10483.3125
$ws.Range("L138").Value = 23839.9761
$ws.Range("M138").Value = -5343.3125
$ws.Range("N138").Value = -34119.9761

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1967.2858
$ws.Range("I63").Value = 2158.8333
$ws.Range("K63").Value = 2158.8333
$ws.Range("M63").Value = -1472.8333

$ws.Range("H66").Value = 1967.2858
$ws.Range("I66").Value = 2158.8333
$ws.Range("K66").Value = 10794.1665
$ws.Range("M66").Value = -7362.166499999999

$ws.Range("H132").Value = 38469800
$ws.Range("I132").Value = 3925.375
$ws.Range("J132").Value = 100015190
$ws.Range("K132").Value = 11776.125
$ws.Range("L132").Value = 300045570
$ws.Range("M132").Value = -9246.125
$ws.Range("N132").Value = -300050630

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H86").Value = 3338
$ws.Range("I86").Value = 1821
$ws.Range("J86").Value = 5866.3335
$ws.Range("K86").Value = 1821
$ws.Range("L86").Value = 5866.3335
$ws.Range("M86").Value = -698
$ws.Range("N86").Value = -8112.3335

$ws.Range("H89").Value = 3338
$ws.Range("I89").Value = 1821
$ws.Range("J89").Value = 5866.3335
$ws.Range("K89").Value = 9105
$ws.Range("L89").Value = 29331.6675
$ws.Range("M89").Value = -3489
$ws.Range("N89").Value = -40563.6675

$ws.Range("H94").Value = 25343.555
$ws.Range("I94").Value = 53799.75
$ws.Range("K94").Value = 53799.75
$ws.Range("M94").Value = -53348.75

$ws.Range("H99").Value = 1899.3334
$ws.Range("I99").Value = 1899.3334
$ws.Range("K99").Value = 1899.3334
$ws.Range("M99").Value = -401.3334

$ws.Range("H134").Value = 50001516
$ws.Range("I134").Value = 50001516
$ws.Range("K134").Value = 150004548
$ws.Range("M134").Value = -150002013

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5010197.5
$ws.Range("I4").Value = 10998.4
$ws.Range("J4").Value = 10009397
$ws.Range("K4").Value = 10998.4
$ws.Range("L4").Value = 10009397
$ws.Range("M4").Value = -10886.4
$ws.Range("N4").Value = -10009621

$ws.Range("H58").Value = 3058.1428
$ws.Range("I58").Value = 3191.4
$ws.Range("J58").Value = 2725
$ws.Range("K58").Value = 3191.4
$ws.Range("L58").Value = 2725
$ws.Range("M58").Value = -2988.4
$ws.Range("N58").Value = -3131

$ws.Range("H62").Value = 20445.445
$ws.Range("I62").Value = 17501.25
$ws.Range("K62").Value = 17501.25
$ws.Range("M62").Value = -16877.25

$ws.Range("H65").Value = 20445.445
$ws.Range("I65").Value = 17501.25
$ws.Range("K65").Value = 87506.25
$ws.Range("M65").Value = -84386.25

$ws.Range("H122").Value = 52099.2
$ws.Range("I122").Value = 2498
$ws.Range("K122").Value = 7494
$ws.Range("M122").Value = -5044

$ws.Range("H136").Value = 3058.1428
$ws.Range("I136").Value = 3191.4
$ws.Range("J136").Value = 2725
$ws.Range("K136").Value = 9574.200000000001
$ws.Range("L136").Value = 8175
$ws.Range("M136").Value = -7024.200000000001
$ws.Range("N136").Value = -13275

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2035.3636
$ws.Range("I5").Value = 1232.8334
$ws.Range("J5").Value = 2998.4
$ws.Range("K5").Value = 3698.5002
$ws.Range("L5").Value = 8995.200000000001
$ws.Range("M5").Value = -3586.5002
$ws.Range("N5").Value = -9219.200000000001

$ws.Range("H15").Value = 114.63636
$ws.Range("I15").Value = 114
$ws.Range("J15").Value = 115
$ws.Range("K15").Value = 342
$ws.Range("L15").Value = 345
$ws.Range("M15").Value = -202
$ws.Range("N15").Value = -625

$ws.Range("H33").Value = 284.44446
$ws.Range("J33").Value = 379
$ws.Range("L33").Value = 2274
$ws.Range("N33").Value = -2840

$ws.Range("H38").Value = 137.35715
$ws.Range("I38").Value = 16.571428
$ws.Range("J38").Value = 258.14285
$ws.Range("K38").Value = 49.71428400000001
$ws.Range("L38").Value = 774.4285500000001
$ws.Range("M38").Value = 297.285716
$ws.Range("N38").Value = -1468.42855

$ws.Range("H68").Value = 3059.8
$ws.Range("J68").Value = 3059.8
$ws.Range("L68").Value = 9179.400000000001
$ws.Range("N68").Value = -10801.4

$ws.Range("H71").Value = 3059.8
$ws.Range("J71").Value = 3059.8
$ws.Range("L71").Value = 27538.2
$ws.Range("N71").Value = -35650.2

$ws.Range("H107").Value = 1876.7727
$ws.Range("J107").Value = 1952.0952
$ws.Range("L107").Value = 5856.2856
$ws.Range("N107").Value = -9696.285599999999

$ws.Range("H122").Value = 62912.75
$ws.Range("I122").Value = 532.4
$ws.Range("J122").Value = 166880
$ws.Range("K122").Value = 4791.599999999999
$ws.Range("L122").Value = 1501920
$ws.Range("M122").Value = -2341.599999999999
$ws.Range("N122").Value = -1506820

$ws.Range("H132").Value = 1779.1
$ws.Range("I132").Value = 1124
$ws.Range("J132").Value = 2215.8333
$ws.Range("K132").Value = 10116
$ws.Range("L132").Value = 19942.4997
$ws.Range("M132").Value = -7586
$ws.Range("N132").Value = -25002.4997

$ws.Range("H135").Value = 2035.3636
$ws.Range("I135").Value = 1232.8334
$ws.Range("J135").Value = 2998.4
$ws.Range("K135").Value = 11095.5006
$ws.Range("L135").Value = 26985.6
$ws.Range("M135").Value = -8560.500599999999
$ws.Range("N135").Value = -32055.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2117.2632
$ws.Range("I102").Value = 1738.2
$ws.Range("K102").Value = 1738.2
$ws.Range("M102").Value = -116.2

$ws.Range("H113").Value = 2901.5557
$ws.Range("I113").Value = 2987.5
$ws.Range("J113").Value = 2729.6667
$ws.Range("K113").Value = 2987.5
$ws.Range("L113").Value = 2729.6667
$ws.Range("M113").Value = -817.5
$ws.Range("N113").Value = -7069.6667

$ws.Range("H132").Value = 3621.7693
$ws.Range("I132").Value = 3525.9048
$ws.Range("J132").Value = 4024.4
$ws.Range("K132").Value = 10577.7144
$ws.Range("L132").Value = 12073.2
$ws.Range("M132").Value = -8047.714399999999
$ws.Range("N132").Value = -17133.2

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7000
$ws.Range("J122").Value = 10000
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -34900

$ws.Range("H132").Value = 4390
$ws.Range("I132").Value = 4380
$ws.Range("K132").Value = 13140
$ws.Range("M132").Value = -10610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 50000000
$ws.Range("J5").Value = 50000000
$ws.Range("L5").Value = 50000000
$ws.Range("N5").Value = -50000224

$ws.Range("H74").Value = 17699
$ws.Range("J74").Value = 18315.5
$ws.Range("L74").Value = 18315.5
$ws.Range("N74").Value = -20187.5

$ws.Range("H77").Value = 17699
$ws.Range("J77").Value = 18315.5
$ws.Range("L77").Value = 54946.5
$ws.Range("N77").Value = -64306.5

$ws.Range("H107").Value = 6886.3335
$ws.Range("I107").Value = 7961.5713
$ws.Range("K107").Value = 23884.7139
$ws.Range("M107").Value = -21964.7139

$ws.Range("H126").Value = 3404.389
$ws.Range("I126").Value = 3165
$ws.Range("J126").Value = 3883.1667
$ws.Range("K126").Value = 9495
$ws.Range("L126").Value = 11649.5001
$ws.Range("M126").Value = -7025
$ws.Range("N126").Value = -16589.5001

$ws.Range("H132").Value = 3445.4075
$ws.Range("J132").Value = 4556.8335
$ws.Range("L132").Value = 13670.5005
$ws.Range("N132").Value = -18730.5005
